# Regenerate save_data: recalculated "K" values (column G) for each row.
# (commit: regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 4
    4  = 1
    5  = 3
    6  = 0
    7  = 1
    9  = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 2
    14 = 1
    15 = 0
    16 = 2
    17 = 2
    18 = 0
    19 = 0
    20 = 0
    21 = 4
    22 = 0
    23 = 2
    24 = 1
    25 = 1
    26 = 2
    27 = 3
    28 = 0
    29 = 1
    30 = 2
    31 = 0
    32 = 0
    33 = 1
    34 = 1
    35 = 2
    36 = 0
    37 = 2
    38 = 2
    39 = 0
    40 = 1
    41 = 2
    42 = 2
    43 = 1
    44 = 2
    45 = 3
    46 = 3
    47 = 3
    48 = 3
    49 = 3
    50 = 1
    51 = 3
    52 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
